# CNNO-2010 - Update HEADERS for file output
#
# Renames several "CBC ..." column headers on the "Data" sheet to more
# generic/neutral names, widens the two "external id" columns (dropping
# their auto/bestFit width in favour of explicit custom widths) and moves
# the active selection from the previously-scrolled last column (U1) back
# to G1.
#
# Header changes (row 1):
#   F1: "CBC Subscription ID" -> "External Subscription ID"
#   G1: "CBC Customer ID"     -> "External Customer ID"
#   M1: "CBC Status"          -> "Status"
#   P1: "CBC Creation Date"   -> "Creation Date"
#   S1: "CBC Licenses"        -> "Licenses"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header labels. Order mirrors the shared-string table layout
# produced by the original edit (M, P, S first, then F, G).
$ws.Range("M1").Value = "Status"
$ws.Range("P1").Value = "Creation Date"
$ws.Range("S1").Value = "Licenses"
$ws.Range("F1").Value = "External Subscription ID"
$ws.Range("G1").Value = "External Customer ID"

# Columns F and G grow to fit the longer, renamed headers - set explicit
# custom widths (no longer auto/best-fit).
$ws.Columns.Item(6).ColumnWidth = 37.584
$ws.Columns.Item(7).ColumnWidth = 54.584

# Move the selection/active cell to G1 (previously the view was scrolled
# all the way to U1).
$ws.Range("G1").Select()
